$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new row is inserted for the "artist" fk entry (was row 7, now row 6),
# the "album" fk entry moves up from row 8 to row 7, and the old row 8
# text cells (artist/idArtist data that used to live there) are cleared,
# leaving only the ImageTrack table header (I8:K8) in place.

# New row 6: artist / idArtist / fk  (copy style from existing row 7 fk cell)
$ws.Range("A6").Value = "artist"
$ws.Range("B6").Value = "idArtist"
$ws.Range("C6").Value = "fk"
$ws.Range("C7").Copy()
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats

# Row 7 becomes: album / idAlbum / fk
$ws.Range("A7").Value = "album"
$ws.Range("B7").Value = "idAlbum"
$ws.Range("C7").Value = "fk"

# Row 8: clear the old album data that used to sit in A8:C8
$ws.Range("A8:C8").Clear()

# Update the active selection to match the saved view state
$ws.Range("K11").Select()
